# Applies the "Suppression des references et correctifs" edit to
# StructureDefinition-Author.xlsx:
#  - Metadata!B8 (Date) -> new timestamp
#  - Metadata!B12 (Description) -> markdown-bold the sub-attribute names in
#    the final paragraph
#  - Elements sheet row for "Author.person":
#      ID/Path/Base Path "Author.person" -> "Author.person[x]"
#      Type(s) "Reference(...)" -> flattened list of profile URLs (no wrapper)
#      Short/Definition "author" -> "Author"
#  - Best-effort re-fit of the Type(s) column width

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- 1. Date -----------------------------------------------------------
$wsMeta.Range("B8").Value = "2025-05-05T11:54:16+00:00"

# --- 2. Description: bold the sub-attribute names -----------------------
# This text is shared by Metadata!B12 and Elements!M2 (the "Author" row's
# Definition) in the source workbook, so both must be updated together.
$descCell = $wsMeta.Range("B12")
$oldDesc = $descCell.Value()
$newDesc = $oldDesc.Replace(
    "author est un ensemble constitu" + [char]0x00E9 + " des sous-attributs authorInstitution , authorPerson, authorRole et authorSpecialty et ne porte pas de valeur par lui-m" + [char]0x00EA + "me. ",
    "**author** est un ensemble constitu" + [char]0x00E9 + " des sous-attributs **authorInstitution** , **authorPerson**, **authorRole** et **authorSpecialty** et ne porte pas de valeur par lui-m" + [char]0x00EA + "me. "
)
$descCell.Value = $newDesc
$wsElem.Range("M2").Value = $newDesc

# --- 3. Elements!row 4 (Author.person -> Author.person[x]) --------------
$wsElem.Range("A4").Value = "Author.person[x]"
$wsElem.Range("B4").Value = "Author.person[x]"
$wsElem.Range("AF4").Value = "Author.person[x]"

$wsElem.Range("K4").Value = "https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/ActorPS`nhttps://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/ActorPatienthttps://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/ActorSNRhttps://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/ActorSystem"

$wsElem.Range("L4").Value = "Author"
$wsElem.Range("M4").Value = "Author"

# --- 4. Best-effort column re-fit ---------------------------------------
$wsElem.Columns.Item(11).ColumnWidth = 189.19140625
